# Adaptations Data Model for Legal info on Resources
# Adds an "Authorship Resource" column (N) to Table1 / Sheet1, populating
# the header and every data row with the resource authorship note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the Excel Table (ListObject) from M52 to N52, 13 -> 14 columns ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N52"))

# --- Header cell (keeps the existing header style s="20") ---
$ws.Range("N1").Value = "Authorship Resource"

# --- Data cells N2:N52 get the shared authorship note ---
$body = $ws.Range("N2:N52")
$body.Value = "Daniela Subotic, Noémi Villars-Amberg"

# Make sure every cell in the new column carries the same look as the rest
# of the data rows (font Arial 14, vertical-top aligned) even for the rows
# that previously had no row-level default formatting (rows 48-52).
$body.Font.Name = "Arial"
$body.Font.Size = 14
$body.VerticalAlignment = -4160

# --- Leave the selection on the newly-added column, matching the saved view ---
$ws.Range("N2:N52").Select()
